$wb = $excel.ActiveWorkbook

# --- Status text update: "Ready for handoff" -> "In Translation" ---
# This text is shared by the per-language status cells on the Overview sheet
# (columns E/F, one per locale) and the "Status" column (C) on each of the
# locale detail sheets. Update every occurrence so they continue to share
# the same underlying text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the per-locale status columns ---
# Overview sheet: columns E and F (the zh-cn / de-de status columns)
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511

# Locale detail sheets: column C (the Status column)
$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511
$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
